$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("000221", "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. baton length of 50. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 12-Mar-2023 17:19:36"),
    @("000222", "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. baton length of 50. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 12-Mar-2023 17:20:31"),
    @("000223", "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. baton length of 50. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 12-Mar-2023 17:24:16"),
    @("000224", "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. baton length of 50. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 12-Mar-2023 17:30:13"),
    @("000225", "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. baton length of 50. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 12-Mar-2023 17:30:59")
)

$startRow = 222
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i

    # Column A: GUID-like code (e.g. "000221"). Force text storage so the
    # leading zeros survive, then strip the formatting override so the
    # cell ends up with the default (General) style, matching the rest
    # of the sheet.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $newRows[$i][0]
    $cellA.ClearFormats()

    # Column B: free-text details string - stored as text naturally.
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
